$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move "Rule Name" label from A8 to A5
$ws.Range("A8").ClearContents()
$ws.Range("A5").Value = "Rule Name"

# Update the active selection as recorded in the saved file (B12)
$ws.Range("B12").Select()
